# Change the table style used by the table on slide 6 ("SOURCES OF FINANCE")
# from the deck's embedded "Table_0" style to the built-in
# "No Style, No Grid" table style ({226E0145-8CA2-41F2-9DEA-FD402EC84C5D}).
#
# Table styles can't be reassigned by setting a property directly - PowerPoint
# requires calling Table.ApplyStyle("{GUID}") with the brace-wrapped StyleId.

$p = $ppt.ActivePresentation

$slide = $p.Slides.Item(6)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{226E0145-8CA2-41F2-9DEA-FD402EC84C5D}")
    }
}
